$d = $word.ActiveDocument

# "Program Time" -> "Programming Time", but expressed (per the target
# OOXML) as a 3-run split of the single original run:
#   "Program" + "ming" + " Time"
# instead of a plain single-run text replace. We rebuild each matching
# paragraph's XML (keeping its original attributes / <w:pPr> untouched)
# with the three runs and push it back in with Range.InsertXML, which
# replaces only the targeted range's contents.

$oldText = "Program Time"
$partName = "/word/document.xml"

$searchRange = $d.Content
while ($searchRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {

    $hitStart = $searchRange.Start
    $hitEnd = $searchRange.End

    # Pull the full package XML fresh each time so offsets reflect any
    # earlier edits in this loop.
    $full = $d.Content.WordOpenXML

    $idx = $full.IndexOf($oldText)
    while ($idx -ge 0) {
        $prefix = $full.Substring(0, $idx)
        $pOpenIdx = $prefix.LastIndexOf("<w:p ")
        $pOpenIdxBare = $prefix.LastIndexOf("<w:p>")
        if ($pOpenIdxBare -gt $pOpenIdx) { $pOpenIdx = $pOpenIdxBare }

        $pCloseIdx = $full.IndexOf(">", $pOpenIdx)
        $openTag = $full.Substring($pOpenIdx, $pCloseIdx - $pOpenIdx + 1)

        $afterOpenTag = $pCloseIdx + 1
        $rest = $full.Substring($afterOpenTag)

        $pPr = ""
        if ($rest.StartsWith("<w:pPr>")) {
            $pPrEnd = $rest.IndexOf("</w:pPr>") + 8
            $pPr = $rest.Substring(0, $pPrEnd)
        }

        # Only accept this paragraph if "Program Time" is still inside it
        # (i.e. no intervening </w:p>) - otherwise keep scanning forward.
        $closeParaIdx = $rest.IndexOf("</w:p>")
        $oldTextIdxInRest = $rest.IndexOf($oldText)
        if ($oldTextIdxInRest -ge 0 -and ($closeParaIdx -lt 0 -or $oldTextIdxInRest -lt $closeParaIdx)) {
            break
        }

        $idx = $full.IndexOf($oldText, $idx + 1)
    }

    $newParaInner = $pPr + "<w:r><w:t>Program</w:t></w:r><w:r><w:t>ming</w:t></w:r><w:r><w:t xml:space=`"preserve`"> Time</w:t></w:r>"
    $newParaXml = $openTag + $newParaInner + "</w:p>"

    $xmlPackage = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="' + $partName + '" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $targetRange = $d.Range($hitStart, $hitEnd)
    $targetRange.InsertXML($xmlPackage)

    # Continue searching after this hit.
    $searchRange = $d.Range($hitStart, $d.Content.End)
}

Write-Output "done"
